# F04 Froze Encoder 1234
# Update the "LJ Speech" worksheet values to reflect the latest ASR results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B text corrections (mis-recognized / corrected words)
$ws.Range("B8").Value  = "<upward>"
$ws.Range("B10").Value = "<see>"
$ws.Range("B12").Value = "<there>"
$ws.Range("B13").Value = "<this>"
$ws.Range("B17").Value = "<him>"
$ws.Range("B18").Value = "<then>"

# Column C numeric score updates
$ws.Range("C2").Value  = 26
$ws.Range("C3").Value  = 33
$ws.Range("C4").Value  = 26
$ws.Range("C5").Value  = 34
$ws.Range("C6").Value  = 22
$ws.Range("C7").Value  = 18
$ws.Range("C8").Value  = 23
$ws.Range("C9").Value  = 27
$ws.Range("C10").Value = 30
$ws.Range("C11").Value = 27
$ws.Range("C12").Value = 39
$ws.Range("C13").Value = 29
$ws.Range("C14").Value = 21
$ws.Range("C15").Value = 24
$ws.Range("C16").Value = 29
$ws.Range("C17").Value = 24
